$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# --- "About" sheet updates ---
$aboutWs = $wb.Worksheets.Item("About")

$aboutWs.Range("A2").Value = "Version: $newVersion"

$aboutWs.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for KWK Borynia-Zofi" + [char]0xF3 + "wka-Bzie Coal Mine, Poland, M1287, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet updates ---
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 24; $row++) {
    $cell = $dataWs.Cells.Item($row, 19)  # column S = 19 (build_version)
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
